# Update the "Estado de Cuenta" worksheet with the new period (2508) and
# updated overdue values (VALOR MORA) for the two worker rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora for the second worker row (E17) changes from 2506 to 2508
$ws.Range("E17").Value = "2508"

# Valor Mora for both worker rows (G16, G17) is updated to the new amount
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
